$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2429577464788732
$ws.Range("C2").Value = 0.4507042253521127
$ws.Range("J2").Value = 0.01408450704225352
$ws.Range("P2").Value = 0.1866197183098592
$ws.Range("S2").Value = 0.1056338028169014
$ws.Range("C3").Value = 0.007874015748031496
$ws.Range("J3").Value = 0.03149606299212598
$ws.Range("P3").Value = 0.8031496062992126
$ws.Range("S3").Value = 0.1574803149606299
$ws.Range("J4").Value = 0.04651162790697674
$ws.Range("P4").Value = 0.5813953488372093
$ws.Range("S4").Value = 0.3720930232558139
$ws.Range("B6").Value = 0.05416666666666667
$ws.Range("D6").Value = 0.008333333333333333
$ws.Range("F6").Value = 0.07083333333333333
$ws.Range("J6").Value = 0.2083333333333333
$ws.Range("O6").Value = 0.008333333333333333
$ws.Range("Q6").Value = 0.2125
$ws.Range("R6").Value = 0.07083333333333333
$ws.Range("S6").Value = 0.3666666666666666
$ws.Range("B7").Value = 0.1022222222222222
$ws.Range("D7").Value = 0.01333333333333333
$ws.Range("F7").Value = 0.07111111111111111
$ws.Range("J7").Value = 0.12
$ws.Range("O7").Value = 0.01777777777777778
$ws.Range("Q7").Value = 0.1688888888888889
$ws.Range("R7").Value = 0.09777777777777778
$ws.Range("S7").Value = 0.4088888888888889
$ws.Range("B8").Value = 0.08806262230919765
$ws.Range("D8").Value = 0.01956947162426614
$ws.Range("E8").Value = 0.001956947162426614
$ws.Range("F8").Value = 0.09001956947162426
$ws.Range("J8").Value = 0.1017612524461839
$ws.Range("O8").Value = 0.01956947162426614
$ws.Range("Q8").Value = 0.1741682974559687
$ws.Range("R8").Value = 0.0821917808219178
$ws.Range("S8").Value = 0.4227005870841487
$ws.Range("B9").Value = 0.1049723756906077
$ws.Range("F9").Value = 0.04972375690607735
$ws.Range("J9").Value = 0.09392265193370165
$ws.Range("O9").Value = 0.02209944751381215
$ws.Range("Q9").Value = 0.1767955801104972
$ws.Range("R9").Value = 0.08287292817679558
$ws.Range("S9").Value = 0.4696132596685083
$ws.Range("B10").Value = 0.08814589665653495
$ws.Range("D10").Value = 0.02203647416413374
$ws.Range("E10").Value = 0.0007598784194528875
$ws.Range("F10").Value = 0.06990881458966565
$ws.Range("J10").Value = 0.1018237082066869
$ws.Range("O10").Value = 0.01595744680851064
$ws.Range("Q10").Value = 0.222644376899696
$ws.Range("R10").Value = 0.0858662613981763
$ws.Range("S10").Value = 0.3928571428571428
$ws.Range("G11").Value = 0.1560846560846561
$ws.Range("J11").Value = 0.08465608465608465
$ws.Range("K11").Value = 0.1984126984126984
$ws.Range("L11").Value = 0.5396825396825397
$ws.Range("S11").Value = 0.02116402116402116
$ws.Range("G12").Value = 0.6807511737089202
$ws.Range("J12").Value = 0.215962441314554
$ws.Range("K12").Value = 0.0187793427230047
$ws.Range("L12").Value = 0.04225352112676056
$ws.Range("S12").Value = 0.04225352112676056
$ws.Range("G13").Value = 0.6444444444444445
$ws.Range("J13").Value = 0.2222222222222222
$ws.Range("S13").Value = 0.1333333333333333
$ws.Range("H15").Value = 0.1617021276595745
$ws.Range("I15").Value = 0.07659574468085106
$ws.Range("J15").Value = 0.3404255319148936
$ws.Range("K15").Value = 0.08936170212765958
$ws.Range("M15").Value = 0.01276595744680851
$ws.Range("O15").Value = 0.02978723404255319
$ws.Range("S15").Value = 0.2893617021276595
$ws.Range("F16").Value = 0.02285714285714286
$ws.Range("H16").Value = 0.1314285714285714
$ws.Range("I16").Value = 0.1028571428571429
$ws.Range("J16").Value = 0.3314285714285714
$ws.Range("K16").Value = 0.1542857142857143
$ws.Range("M16").Value = 0.02285714285714286
$ws.Range("O16").Value = 0.04
$ws.Range("S16").Value = 0.1942857142857143
$ws.Range("F17").Value = 0.01414141414141414
$ws.Range("H17").Value = 0.195959595959596
$ws.Range("I17").Value = 0.07676767676767676
$ws.Range("J17").Value = 0.404040404040404
$ws.Range("K17").Value = 0.1171717171717172
$ws.Range("M17").Value = 0.01616161616161616
$ws.Range("O17").Value = 0.04848484848484848
$ws.Range("S17").Value = 0.1272727272727273
$ws.Range("F18").Value = 0.01428571428571429
$ws.Range("H18").Value = 0.1761904761904762
$ws.Range("I18").Value = 0.0761904761904762
$ws.Range("J18").Value = 0.4047619047619048
$ws.Range("K18").Value = 0.08571428571428572
$ws.Range("M18").Value = 0.009523809523809525
$ws.Range("O18").Value = 0.08095238095238096
$ws.Range("S18").Value = 0.1523809523809524
$ws.Range("F19").Value = 0.01889433170048985
$ws.Range("H19").Value = 0.2232330300909727
$ws.Range("I19").Value = 0.06508047585724283
$ws.Range("J19").Value = 0.3785864240727781
$ws.Range("K19").Value = 0.1147655703289013
$ws.Range("M19").Value = 0.01819454163750875
$ws.Range("O19").Value = 0.07207837648705388
$ws.Range("S19").Value = 0.1091672498250525

Write-Host "Applied team specific time data updates"
